$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.038426160812378
$ws.Range("B1").Value = 4.346137046813965
$ws.Range("C1").Value = 3.435123443603516
$ws.Range("D1").Value = 2.218516111373901
$ws.Range("E1").Value = 2.037258625030518
